$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.708.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "'1.607.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.65%  "

$ws.Range("D5").Value = "'212.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("E6").Value = "  +1.22%  "

$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("D8").Value = "'28.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.11%  "

$ws.Range("E9").Value = "  +1.31%  "

$ws.Range("E10").Value = "  +1.46%  "

$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").Value = "'1.837.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.36%  "

$ws.Range("D13").Value = "'1.583.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("E14").Value = "  +3.94%  "

$ws.Range("D15").Value = "'29.717.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "'64.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "'241.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").Value = "'7.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.59%  "

$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("E21").Value = "  +0.80%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "'9.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "

$ws.Range("D24").Value = "'2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("D25").Value = "'155.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("E26").Value = "  +1.64%  "

$ws.Range("D27").Value = "'0.109"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("D28").Value = "'6.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.19%  "

$ws.Range("E29").Value = "  +0.61%  "

$ws.Range("E30").Value = "  +1.92%  "

$ws.Range("E31").Value = "  +1.29%  "

$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("E33").Value = "  +2.60%  "

$ws.Range("D34").Value = "'1.430.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").Value = "'2.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.62%  "

$ws.Range("E36").Value = "  +3.12%  "

$ws.Range("D37").Value = "'1.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("B40").Value = "BitcoinSV"
$ws.Range("C40").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D40").Value = "'57.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.63%  "

$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'0.549"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.69%  "

$ws.Range("D42").Value = "'0.0499"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.23%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.819"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.84%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").Value = "'66.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.26%  "

$ws.Range("D47").Value = "'0.983"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.27%  "

$ws.Range("D48").Value = "'5.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.56%  "

$ws.Range("D49").Value = "'1.746.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "

$ws.Range("D50").Value = "'86.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("E51").Value = "  +3.71%  "
